# Task: "...Upon visiting the website, it was very clear..." becomes
# "...Upon visiting the website and joining the server, it was very clear..."
# and Word's auto-managed "_GoBack" bookmark (which tracks the location of the
# most recent edit) moves from its old spot (an empty paragraph further down
# the document) to sit right after the newly typed text.
#
# Note: this engine re-coalesces adjacent same-formatted runs whenever a run's
# text is spliced in-place, so a naive Range.Text/InsertBefore/InsertAfter at
# the target offset would merge across run boundaries that must stay distinct
# (per the target XML). Wrapping a (possibly zero-length) sub-range with a
# throwaway Bookmarks.Add/Delete pins a hard run-split at its edges without
# touching the text -- and that split survives deleting the bookmark again --
# so we use that trick to keep "l. Upon visiting the website", the newly
# inserted " and joining the server", and the remainder of the sentence as
# three separate runs, exactly like Word itself would leave them after a
# mid-sentence insertion.

$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark (currently sat alone in
#     an empty paragraph near the end of the document). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate "l. Upon visiting the website" and pin a hard run-split
#     between it and the run preceding it, without altering any text. ---
$rWrap = $d.Content
$rWrap.Find.Execute("l. Upon visiting the website")
$wrapEnd = $rWrap.End
$d.Bookmarks.Add("__tmp_split_before", $rWrap) | Out-Null

# --- Step 3: insert a unique placeholder right after "website" (i.e. right
#     before the comma) that will become the new run's text. ---
$rIns = $d.Range($wrapEnd, $wrapEnd)
$rIns.InsertBefore("@@NEWTEXT@@")

# --- Step 4: find the placeholder, wrap it with a bookmark (pinning splits on
#     both sides of it), then swap its text for the real sentence addition. ---
$rPh = $d.Content
$rPh.Find.Execute("@@NEWTEXT@@")
$bmPh = $d.Bookmarks.Add("__tmp_placeholder", $rPh)
$bmPh.Range.Text = " and joining the server"

# --- Step 5: re-create "_GoBack" as a zero-length marker right after the
#     newly inserted text, matching Word's "last edit position" bookmark. ---
$bmPh2 = $d.Bookmarks.Item("__tmp_placeholder")
$afterPos = $bmPh2.End
$rGoBack = $d.Range($afterPos, $afterPos)
$d.Bookmarks.Add("_GoBack", $rGoBack) | Out-Null

# --- Step 6: clean up the throwaway bookmarks that only existed to pin run
#     boundaries; the splits they created remain even after deletion. ---
$d.Bookmarks.Item("__tmp_split_before").Delete()
$d.Bookmarks.Item("__tmp_placeholder").Delete()
